$wb = $excel.ActiveWorkbook

# Updates to the "想去人数" (F column) figures, applied identically to
# both the "展览" and "全部类型" worksheets.
$updates = @{
    "F2"  = 8460
    "F3"  = 8130
    "F9"  = 143
    "F10" = 196
    "F13" = 195
    "F14" = 4212
    "F16" = 73
    "F19" = 148
    "F20" = 109
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
